# Generate Report for Handoff
# A new handoff run produced a new source-file GUID + commit hash, refreshed the
# handoff timestamps for both locales, and reset the (not-yet-done) handback
# fields back to "empty" / the epoch sentinel date.

$wb = $excel.ActiveWorkbook

$oldGuid = "17d42d58-5e2d-4784-9aef-f69d7b3e93c1"
$newGuid = "c127dd0e-7e76-4186-a688-12796be6f1ae"
$oldHash = "85ece963d56a4f108e92a8f30855354f066ec6eb"
$newHash = "37426f6a6d7e43776475af2d6fcb24e7ac0eeffc"

$newFileName  = "$newGuid.md"
$newPathName  = "e2e\$newGuid.md"
$newHoDate    = "2016-09-07 07:24:40"
$neverDate    = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("A2").Value = $newFileName
$overview.Range("B2").Value = $newPathName
$overview.Range("G2").Value = $newHoDate

# Keep the hyperlink's displayed text in sync with the new path/name.
$overview.Hyperlinks.Item(1).TextToDisplay = $newPathName

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("A2").Value = $newFileName
$zhcn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$zhcn.Range("H2").Value = "2016-09-07 07:24:34"
$zhcn.Range("I2").Value = ""
$zhcn.Range("J2").Value = ""
$zhcn.Range("K2").Value = $neverDate

$zhcn.Columns.Item(9).ColumnWidth = 18.6506053379604
$zhcn.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("A2").Value = $newFileName
$dede.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$dede.Range("H2").Value = $newHoDate
$dede.Range("I2").Value = ""
$dede.Range("J2").Value = ""
$dede.Range("K2").Value = $neverDate

$dede.Columns.Item(9).ColumnWidth = 18.6506053379604
$dede.Columns.Item(10).ColumnWidth = 21.7054770333426
